$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($totalSheet)
$ws = $wb.Worksheets.Item("2021-Q4 (2)")
$ws.Name = "2022-Q1"

# Trim the copied sheet from 90 data rows down to the 40 we need.
$ws.Rows("42:91").Delete()

# Helper: write $val into $cell as TEXT (avoids Excel auto-numeric coercion)
# by round-tripping it through a quote-prefixed scratch cell + values-only paste.
$scratch = $ws.Cells.Item(200, 50)
function Set-TextValue($cell, $val) {
    $scratch.Value = "'" + $val
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextValue $ws.Cells.Item(2,2) '010488'
Set-TextValue $ws.Cells.Item(2,3) '鹏华优选成长混合A'
Set-TextValue $ws.Cells.Item(2,4) '34.99'
Set-TextValue $ws.Cells.Item(2,5) '69.88'
Set-TextValue $ws.Cells.Item(2,6) '3.10'
Set-TextValue $ws.Cells.Item(2,7) '1.0847'
$ws.Cells.Item(2,8).Value = 5

Set-TextValue $ws.Cells.Item(3,2) '001230'
Set-TextValue $ws.Cells.Item(3,3) '鹏华医药科技股票'
Set-TextValue $ws.Cells.Item(3,4) '16.65'
Set-TextValue $ws.Cells.Item(3,5) '81.88'
Set-TextValue $ws.Cells.Item(3,6) '5.12'
Set-TextValue $ws.Cells.Item(3,7) '0.8525'
$ws.Cells.Item(3,8).Value = 5

Set-TextValue $ws.Cells.Item(4,2) '501011'
Set-TextValue $ws.Cells.Item(4,3) '汇添富中证中药指数（LOF）A'
Set-TextValue $ws.Cells.Item(4,4) '24.13'
Set-TextValue $ws.Cells.Item(4,5) '94.57'
Set-TextValue $ws.Cells.Item(4,6) '3.18'
Set-TextValue $ws.Cells.Item(4,7) '0.7673'
$ws.Cells.Item(4,8).Value = 9

Set-TextValue $ws.Cells.Item(5,2) '011568'
Set-TextValue $ws.Cells.Item(5,3) '鹏华产业升级混合A'
Set-TextValue $ws.Cells.Item(5,4) '21.41'
Set-TextValue $ws.Cells.Item(5,5) '68.03'
Set-TextValue $ws.Cells.Item(5,6) '3.20'
Set-TextValue $ws.Cells.Item(5,7) '0.6851'
$ws.Cells.Item(5,8).Value = 4

Set-TextValue $ws.Cells.Item(6,2) '000242'
Set-TextValue $ws.Cells.Item(6,3) '景顺长城策略精选'
Set-TextValue $ws.Cells.Item(6,4) '20.09'
Set-TextValue $ws.Cells.Item(6,5) '91.24'
Set-TextValue $ws.Cells.Item(6,6) '3.23'
Set-TextValue $ws.Cells.Item(6,7) '0.6489'
$ws.Cells.Item(6,8).Value = 10

Set-TextValue $ws.Cells.Item(7,2) '013414'
Set-TextValue $ws.Cells.Item(7,3) '太平智远三个月定期开放股票'
Set-TextValue $ws.Cells.Item(7,4) '8.69'
Set-TextValue $ws.Cells.Item(7,5) '86.34'
Set-TextValue $ws.Cells.Item(7,6) '4.57'
Set-TextValue $ws.Cells.Item(7,7) '0.3971'
$ws.Cells.Item(7,8).Value = 9

Set-TextValue $ws.Cells.Item(8,2) '010405'
Set-TextValue $ws.Cells.Item(8,3) '惠升医药健康6个月持有期混合'
Set-TextValue $ws.Cells.Item(8,4) '12.18'
Set-TextValue $ws.Cells.Item(8,5) '64.63'
Set-TextValue $ws.Cells.Item(8,6) '3.19'
Set-TextValue $ws.Cells.Item(8,7) '0.3885'
$ws.Cells.Item(8,8).Value = 7

Set-TextValue $ws.Cells.Item(9,2) '160610'
Set-TextValue $ws.Cells.Item(9,3) '鹏华动力增长混合(LOF)'
Set-TextValue $ws.Cells.Item(9,4) '14.23'
Set-TextValue $ws.Cells.Item(9,5) '57.89'
Set-TextValue $ws.Cells.Item(9,6) '2.46'
Set-TextValue $ws.Cells.Item(9,7) '0.3501'
$ws.Cells.Item(9,8).Value = 7

Set-TextValue $ws.Cells.Item(10,2) '003713'
Set-TextValue $ws.Cells.Item(10,3) '英大睿盛灵活配置混合A'
Set-TextValue $ws.Cells.Item(10,4) '5.99'
Set-TextValue $ws.Cells.Item(10,5) '87.42'
Set-TextValue $ws.Cells.Item(10,6) '5.11'
Set-TextValue $ws.Cells.Item(10,7) '0.3061'
$ws.Cells.Item(10,8).Value = 6

Set-TextValue $ws.Cells.Item(11,2) '012093'
Set-TextValue $ws.Cells.Item(11,3) '鹏华创新升级混合型证券投资基金A'
Set-TextValue $ws.Cells.Item(11,4) '6.58'
Set-TextValue $ws.Cells.Item(11,5) '64.31'
Set-TextValue $ws.Cells.Item(11,6) '4.57'
Set-TextValue $ws.Cells.Item(11,7) '0.3007'
$ws.Cells.Item(11,8).Value = 4

Set-TextValue $ws.Cells.Item(12,2) '501012'
Set-TextValue $ws.Cells.Item(12,3) '汇添富中证中药指数（LOF）C'
Set-TextValue $ws.Cells.Item(12,4) '8.91'
Set-TextValue $ws.Cells.Item(12,5) '94.57'
Set-TextValue $ws.Cells.Item(12,6) '3.18'
Set-TextValue $ws.Cells.Item(12,7) '0.2833'
$ws.Cells.Item(12,8).Value = 9

Set-TextValue $ws.Cells.Item(13,2) '000780'
Set-TextValue $ws.Cells.Item(13,3) '鹏华医疗保健股票'
Set-TextValue $ws.Cells.Item(13,4) '7.67'
Set-TextValue $ws.Cells.Item(13,5) '82.80'
Set-TextValue $ws.Cells.Item(13,6) '3.21'
Set-TextValue $ws.Cells.Item(13,7) '0.2462'
$ws.Cells.Item(13,8).Value = 7

Set-TextValue $ws.Cells.Item(14,2) '217001'
Set-TextValue $ws.Cells.Item(14,3) '招商安泰混合'
Set-TextValue $ws.Cells.Item(14,4) '4.22'
Set-TextValue $ws.Cells.Item(14,5) '70.79'
Set-TextValue $ws.Cells.Item(14,6) '4.02'
Set-TextValue $ws.Cells.Item(14,7) '0.1696'
$ws.Cells.Item(14,8).Value = 1

Set-TextValue $ws.Cells.Item(15,2) '160603'
Set-TextValue $ws.Cells.Item(15,3) '鹏华普天收益混合'
Set-TextValue $ws.Cells.Item(15,4) '5.22'
Set-TextValue $ws.Cells.Item(15,5) '65.41'
Set-TextValue $ws.Cells.Item(15,6) '2.93'
Set-TextValue $ws.Cells.Item(15,7) '0.1529'
$ws.Cells.Item(15,8).Value = 6

Set-TextValue $ws.Cells.Item(16,2) '005270'
Set-TextValue $ws.Cells.Item(16,3) '太平改革红利精选灵活配置混合'
Set-TextValue $ws.Cells.Item(16,4) '1.87'
Set-TextValue $ws.Cells.Item(16,5) '88.32'
Set-TextValue $ws.Cells.Item(16,6) '6.57'
Set-TextValue $ws.Cells.Item(16,7) '0.1229'
$ws.Cells.Item(16,8).Value = 6

Set-TextValue $ws.Cells.Item(17,2) '003714'
Set-TextValue $ws.Cells.Item(17,3) '英大睿盛灵活配置混合C'
Set-TextValue $ws.Cells.Item(17,4) '2.40'
Set-TextValue $ws.Cells.Item(17,5) '87.42'
Set-TextValue $ws.Cells.Item(17,6) '5.11'
Set-TextValue $ws.Cells.Item(17,7) '0.1226'
$ws.Cells.Item(17,8).Value = 6

Set-TextValue $ws.Cells.Item(18,2) '012506'
Set-TextValue $ws.Cells.Item(18,3) '东方品质消费一年持有期混合型证券投资基金A'
Set-TextValue $ws.Cells.Item(18,4) '2.51'
Set-TextValue $ws.Cells.Item(18,5) '90.03'
Set-TextValue $ws.Cells.Item(18,6) '3.38'
Set-TextValue $ws.Cells.Item(18,7) '0.0848'
$ws.Cells.Item(18,8).Value = 8

Set-TextValue $ws.Cells.Item(19,2) '002259'
Set-TextValue $ws.Cells.Item(19,3) '鹏华健康环保灵活配置混合'
Set-TextValue $ws.Cells.Item(19,4) '2.19'
Set-TextValue $ws.Cells.Item(19,5) '79.82'
Set-TextValue $ws.Cells.Item(19,6) '3.20'
Set-TextValue $ws.Cells.Item(19,7) '0.0701'
$ws.Cells.Item(19,8).Value = 7

Set-TextValue $ws.Cells.Item(20,2) '011331'
Set-TextValue $ws.Cells.Item(20,3) '鹏华远见成长混合型证券投资基金A'
Set-TextValue $ws.Cells.Item(20,4) '2.29'
Set-TextValue $ws.Cells.Item(20,5) '63.89'
Set-TextValue $ws.Cells.Item(20,6) '3.02'
Set-TextValue $ws.Cells.Item(20,7) '0.0692'
$ws.Cells.Item(20,8).Value = 4

Set-TextValue $ws.Cells.Item(21,2) '400025'
Set-TextValue $ws.Cells.Item(21,3) '东方新兴成长混合'
Set-TextValue $ws.Cells.Item(21,4) '2.19'
Set-TextValue $ws.Cells.Item(21,5) '89.57'
Set-TextValue $ws.Cells.Item(21,6) '3.09'
Set-TextValue $ws.Cells.Item(21,7) '0.0677'
$ws.Cells.Item(21,8).Value = 10

Set-TextValue $ws.Cells.Item(22,2) '005112'
Set-TextValue $ws.Cells.Item(22,3) '银华中证全指医药卫生指数增强'
Set-TextValue $ws.Cells.Item(22,4) '1.42'
Set-TextValue $ws.Cells.Item(22,5) '86.59'
Set-TextValue $ws.Cells.Item(22,6) '4.76'
Set-TextValue $ws.Cells.Item(22,7) '0.0676'
$ws.Cells.Item(22,8).Value = 8

Set-TextValue $ws.Cells.Item(23,2) '010896'
Set-TextValue $ws.Cells.Item(23,3) '太平价值增长股票A'
Set-TextValue $ws.Cells.Item(23,4) '1.18'
Set-TextValue $ws.Cells.Item(23,5) '83.63'
Set-TextValue $ws.Cells.Item(23,6) '5.39'
Set-TextValue $ws.Cells.Item(23,7) '0.0636'
$ws.Cells.Item(23,8).Value = 6

Set-TextValue $ws.Cells.Item(24,2) '519959'
Set-TextValue $ws.Cells.Item(24,3) '长信多利灵活配置混合'
Set-TextValue $ws.Cells.Item(24,4) '1.45'
Set-TextValue $ws.Cells.Item(24,5) '85.11'
Set-TextValue $ws.Cells.Item(24,6) '4.36'
Set-TextValue $ws.Cells.Item(24,7) '0.0632'
$ws.Cells.Item(24,8).Value = 6

Set-TextValue $ws.Cells.Item(25,2) '013488'
Set-TextValue $ws.Cells.Item(25,3) '长信多利灵活配置混合D'
Set-TextValue $ws.Cells.Item(25,4) '1.45'
Set-TextValue $ws.Cells.Item(25,5) '85.11'
Set-TextValue $ws.Cells.Item(25,6) '4.36'
Set-TextValue $ws.Cells.Item(25,7) '0.0632'
$ws.Cells.Item(25,8).Value = 6

Set-TextValue $ws.Cells.Item(26,2) '010897'
Set-TextValue $ws.Cells.Item(26,3) '太平价值增长股票C'
Set-TextValue $ws.Cells.Item(26,4) '1.01'
Set-TextValue $ws.Cells.Item(26,5) '83.63'
Set-TextValue $ws.Cells.Item(26,6) '5.39'
Set-TextValue $ws.Cells.Item(26,7) '0.0544'
$ws.Cells.Item(26,8).Value = 6

Set-TextValue $ws.Cells.Item(27,2) '003446'
Set-TextValue $ws.Cells.Item(27,3) '英大睿鑫灵活配置混合A'
Set-TextValue $ws.Cells.Item(27,4) '0.59'
Set-TextValue $ws.Cells.Item(27,5) '89.46'
Set-TextValue $ws.Cells.Item(27,6) '5.53'
Set-TextValue $ws.Cells.Item(27,7) '0.0326'
$ws.Cells.Item(27,8).Value = 7

Set-TextValue $ws.Cells.Item(28,2) '010489'
Set-TextValue $ws.Cells.Item(28,3) '鹏华优选成长混合C'
Set-TextValue $ws.Cells.Item(28,4) '0.92'
Set-TextValue $ws.Cells.Item(28,5) '69.88'
Set-TextValue $ws.Cells.Item(28,6) '3.10'
Set-TextValue $ws.Cells.Item(28,7) '0.0285'
$ws.Cells.Item(28,8).Value = 5

Set-TextValue $ws.Cells.Item(29,2) '003447'
Set-TextValue $ws.Cells.Item(29,3) '英大睿鑫灵活配置混合C'
Set-TextValue $ws.Cells.Item(29,4) '0.51'
Set-TextValue $ws.Cells.Item(29,5) '89.46'
Set-TextValue $ws.Cells.Item(29,6) '5.53'
Set-TextValue $ws.Cells.Item(29,7) '0.0282'
$ws.Cells.Item(29,8).Value = 7

Set-TextValue $ws.Cells.Item(30,2) '217021'
Set-TextValue $ws.Cells.Item(30,3) '招商优势企业混合'
Set-TextValue $ws.Cells.Item(30,4) '0.36'
Set-TextValue $ws.Cells.Item(30,5) '69.72'
Set-TextValue $ws.Cells.Item(30,6) '4.43'
Set-TextValue $ws.Cells.Item(30,7) '0.0159'
$ws.Cells.Item(30,8).Value = 9

Set-TextValue $ws.Cells.Item(31,2) '012507'
Set-TextValue $ws.Cells.Item(31,3) '东方品质消费一年持有期混合型证券投资基金C'
Set-TextValue $ws.Cells.Item(31,4) '0.42'
Set-TextValue $ws.Cells.Item(31,5) '90.03'
Set-TextValue $ws.Cells.Item(31,6) '3.38'
Set-TextValue $ws.Cells.Item(31,7) '0.0142'
$ws.Cells.Item(31,8).Value = 8

Set-TextValue $ws.Cells.Item(32,2) '519987'
Set-TextValue $ws.Cells.Item(32,3) '长信恒利优势混合'
Set-TextValue $ws.Cells.Item(32,4) '0.22'
Set-TextValue $ws.Cells.Item(32,5) '82.39'
Set-TextValue $ws.Cells.Item(32,6) '5.07'
Set-TextValue $ws.Cells.Item(32,7) '0.0112'
$ws.Cells.Item(32,8).Value = 4

Set-TextValue $ws.Cells.Item(33,2) '011569'
Set-TextValue $ws.Cells.Item(33,3) '鹏华产业升级混合C'
Set-TextValue $ws.Cells.Item(33,4) '0.34'
Set-TextValue $ws.Cells.Item(33,5) '68.03'
Set-TextValue $ws.Cells.Item(33,6) '3.20'
Set-TextValue $ws.Cells.Item(33,7) '0.0109'
$ws.Cells.Item(33,8).Value = 4

Set-TextValue $ws.Cells.Item(34,2) '012094'
Set-TextValue $ws.Cells.Item(34,3) '鹏华创新升级混合型证券投资基金C'
Set-TextValue $ws.Cells.Item(34,4) '0.19'
Set-TextValue $ws.Cells.Item(34,5) '64.31'
Set-TextValue $ws.Cells.Item(34,6) '4.57'
Set-TextValue $ws.Cells.Item(34,7) '0.0087'
$ws.Cells.Item(34,8).Value = 4

Set-TextValue $ws.Cells.Item(35,2) '011332'
Set-TextValue $ws.Cells.Item(35,3) '鹏华远见成长混合型证券投资基金C'
Set-TextValue $ws.Cells.Item(35,4) '0.21'
Set-TextValue $ws.Cells.Item(35,5) '63.89'
Set-TextValue $ws.Cells.Item(35,6) '3.02'
Set-TextValue $ws.Cells.Item(35,7) '0.0063'
$ws.Cells.Item(35,8).Value = 4

Set-TextValue $ws.Cells.Item(36,2) '011548'
Set-TextValue $ws.Cells.Item(36,3) '九泰久慧混合A'
Set-TextValue $ws.Cells.Item(36,4) '0.49'
Set-TextValue $ws.Cells.Item(36,5) '31.08'
Set-TextValue $ws.Cells.Item(36,6) '0.83'
Set-TextValue $ws.Cells.Item(36,7) '0.0041'
$ws.Cells.Item(36,8).Value = 9

Set-TextValue $ws.Cells.Item(37,2) '004917'
Set-TextValue $ws.Cells.Item(37,3) '中银证券祥瑞混合A'
Set-TextValue $ws.Cells.Item(37,4) '0.10'
Set-TextValue $ws.Cells.Item(37,5) '79.01'
Set-TextValue $ws.Cells.Item(37,6) '2.08'
Set-TextValue $ws.Cells.Item(37,7) '0.0021'
$ws.Cells.Item(37,8).Value = 8

Set-TextValue $ws.Cells.Item(38,2) '001608'
Set-TextValue $ws.Cells.Item(38,3) '英大策略优选混合C'
Set-TextValue $ws.Cells.Item(38,4) '0.03'
Set-TextValue $ws.Cells.Item(38,5) '89.86'
Set-TextValue $ws.Cells.Item(38,6) '5.13'
Set-TextValue $ws.Cells.Item(38,7) '0.0015'
$ws.Cells.Item(38,8).Value = 7

Set-TextValue $ws.Cells.Item(39,2) '004918'
Set-TextValue $ws.Cells.Item(39,3) '中银证券祥瑞混合C'
Set-TextValue $ws.Cells.Item(39,4) '0.07'
Set-TextValue $ws.Cells.Item(39,5) '79.01'
Set-TextValue $ws.Cells.Item(39,6) '2.08'
Set-TextValue $ws.Cells.Item(39,7) '0.0015'
$ws.Cells.Item(39,8).Value = 8

Set-TextValue $ws.Cells.Item(40,2) '006195'
Set-TextValue $ws.Cells.Item(40,3) '国金量化多因子股票'
Set-TextValue $ws.Cells.Item(40,4) '0.09'
Set-TextValue $ws.Cells.Item(40,5) '80.71'
Set-TextValue $ws.Cells.Item(40,6) '0.88'
Set-TextValue $ws.Cells.Item(40,7) '0.0008'
$ws.Cells.Item(40,8).Value = 8

Set-TextValue $ws.Cells.Item(41,2) '011549'
Set-TextValue $ws.Cells.Item(41,3) '九泰久慧混合C'
Set-TextValue $ws.Cells.Item(41,4) '0.06'
Set-TextValue $ws.Cells.Item(41,5) '31.08'
Set-TextValue $ws.Cells.Item(41,6) '0.83'
Set-TextValue $ws.Cells.Item(41,7) '0.0005'
$ws.Cells.Item(41,8).Value = 9

$scratch.Clear()

$totalSheet = $wb.Worksheets.Item("总计")

# Extend the styled index-column (A) format down into the new row 7.
$totalSheet.Cells.Item(6,1).Copy()
$totalSheet.Cells.Item(7,1).PasteSpecial(-4122)

# Shift existing data rows 2..6 down to 3..7 (bottom-up, values only).
for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 1
    $totalSheet.Cells.Item($dest,2).Value = $totalSheet.Cells.Item($r,2).Value()
    $totalSheet.Cells.Item($dest,3).Value = $totalSheet.Cells.Item($r,3).Value()
    $totalSheet.Cells.Item($dest,4).Value = $totalSheet.Cells.Item($r,4).Value()
}

# Renumber the 0-based index column for rows 2..7.
for ($r = 2; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r,1).Value = $r - 2
}

# Write the new 2022-Q1 summary into row 2.
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 40
$totalSheet.Cells.Item(2,4).Value = 7.65

